$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column C to make room for the new
# "Group" column; this shifts the existing Amount / Type of Expense /
# Date of Transaction columns from C,D,E to D,E,F (preserving their
# original cell types/styles, since it's a true column insert rather
# than a value-by-value rewrite).
$ws.Range("C1").EntireColumn.Insert()

# The insert pushed the old (empty) trailing cell from F1 into G1;
# clear it so the sheet doesn't pick up a stray 7th column of data.
$ws.Range("G1").Clear()

# New "Group" column header and values
$ws.Cells.Item(1, 3).Value = "Group"
$ws.Cells.Item(2, 3).Value = "ssl2"
$ws.Cells.Item(3, 3).Value = "ssl2"
$ws.Cells.Item(4, 3).Value = "ssl2"
$ws.Cells.Item(5, 3).Value = "N/A"
$ws.Cells.Item(6, 3).Value = "N/A"
$ws.Cells.Item(7, 3).Value = "N/A"
$ws.Cells.Item(8, 3).Value = "N/A"
$ws.Cells.Item(9, 3).Value = "che"
$ws.Cells.Item(10, 3).Value = "ssl work"
$ws.Cells.Item(11, 3).Value = "ssl"

# Column widths: originally A=20 B=20 C=20 D=20 E=30 F=20.
# After inserting the Group column, give it the same width (20) as the
# other columns, and swap the old Type/Date widths so E=20, F=30.
$ws.Columns.Item(3).ColumnWidth = 19.1
$ws.Columns.Item(5).ColumnWidth = 19.1
$ws.Columns.Item(6).ColumnWidth = 29.1
